$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New inventory row for "Frozen Berries"
$ws.Range("A2").Value = "Frozen Berries"
$ws.Range("B2").Value = 1234
$ws.Range("C2").Value = "F"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 562.62
$ws.Range("G2").Value = 35.86
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 100

# Column A needs to widen to fit the new "Frozen Berries" text (bestFit/autofit)
$ws.Columns.Item(1).ColumnWidth = 13
